$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($i=1; $i -le $ws.Shapes.Count; $i++) {
    $shp = $ws.Shapes.Item($i)
    Write-Host $i $shp.Name "Top=" $shp.Top "Left=" $shp.Left "Width=" $shp.Width "Height=" $shp.Height
}
